#
# Adds the new "2023" column (U) to the mortality-rate table, mirroring the
# formatting of the existing "2022" column (T) for each data row, and moves
# the active selection back to B1 (top-left of the visible/frozen area).
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new value for column U (mirrors column T's formatting)
$newValues = @{
    4  = 2023
    5  = 3.3
    6  = 1
    7  = 1.6
    8  = 9.1999999999999993
    9  = 6.1
    10 = 1.5
    11 = 4
    12 = 4.4000000000000004
    13 = 4.7
    14 = 0.5
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("T$row")
    $dstCell = $ws.Range("U$row")

    # Copy column T's formatting (number format, font, borders, alignment, ...)
    # onto the new column U cell for this row.
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)

    # Now write the cell's own value (PasteSpecial above only brought over
    # formatting, not T's value).
    $dstCell.Value = $newValues[$row]
}

$excel.CutCopyMode = $false

# Restore the view/selection to B1, matching the saved workbook state.
$ws.Range("B1").Select()
